# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets: row 2 (F2) 212 -> 213, row 3 (F3) 158 -> 159.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 213
    $ws.Range("F3").Value = 159
}
